$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the resource table (A2:L17) alphabetically by Resource (column A),
# mirroring the "Data > Sort" operation performed when hooking up the new
# resources/hazards into WorldMap_UI.
$sortRange = $ws.Range("A2:L17")
$sortKey   = $ws.Range("A2:A17")

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($sortKey)
$sortObj.SetRange($sortRange)
$sortObj.Header = 2
$sortObj.Apply()

# Leave the same selection state the sheet ended up with after the edit:
# the whole of row 10 selected, with the view scrolled so column E is the
# left-most visible column.
$ws.Range("A10:XFD10").Select()
